# Natmi following Dr Hou advice
# Update LR-pair table: add "ECs" sending/target cluster group and refresh computed stats
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs","Col18a1","Gpc4","ECs",2,0.6666666666666666,10.92859066666667,32.785772,0.2185558471001832,0.2185558471001832,3,1,8.767882,26.303646,0.2073466531282026,0.2073466531282026,95.82059339163467,862.3853405247121,0.04531682341782217,0.04531682341782218),
  @("ECs","Col18a1","Gpc4","FAPs",2,0.6666666666666666,10.92859066666667,32.785772,0.2185558471001832,0.2185558471001832,3,1,20.88867166666667,62.666015,0.4939843121037951,0.4939843121037952,228.2837422153978,2054.55367993858,0.1079631597860462,0.1079631597860462),
  @("ECs","Col18a1","Gpc4","sCs",2,0.6666666666666666,10.92859066666667,32.785772,0.2185558471001832,0.2185558471001832,3,1,12.62954966666667,37.888649,0.2986690347680022,0.2986690347680022,138.0231786113364,1242.208607502028,0.06527586389631478,0.06527586389631479),
  @("FAPs","Col18a1","Gpc4","ECs",3,1,21.275312,63.825936,0.4254751576214852,0.4254751576214852,3,1,8.767882,26.303646,0.2073466531282026,0.2073466531282026,186.539425129184,1678.854826162656,0.08822084992200943,0.08822084992200944),
  @("FAPs","Col18a1","Gpc4","FAPs",3,1,21.275312,63.825936,0.4254751576214852,0.4254751576214852,3,1,20.88867166666667,62.666015,0.4939843121037951,0.4939843121037952,444.4130069738933,3999.71706276504,0.2101780530549032,0.2101780530549032),
  @("FAPs","Col18a1","Gpc4","sCs",3,1,21.275312,63.825936,0.4254751576214852,0.4254751576214852,3,1,12.62954966666667,37.888649,0.2986690347680022,0.2986690347680022,268.6976095778293,2418.278486200464,0.1270762546445726,0.1270762546445726),
  @("sCs","Col18a1","Gpc4","ECs",3,1,17.79975,53.39925,0.3559689952783316,0.3559689952783316,3,1,8.767882,26.303646,0.2073466531282026,0.2073466531282026,156.0661076295,1404.5949686655,0.07380897978837102,0.07380897978837103),
  @("sCs","Col18a1","Gpc4","FAPs",3,1,17.79975,53.39925,0.3559689952783316,0.3559689952783316,3,1,20.88867166666667,62.666015,0.4939843121037951,0.4939843121037952,371.81313349875,3346.31820148875,0.1758430992628457,0.1758430992628458),
  @("sCs","Col18a1","Gpc4","sCs",3,1,17.79975,53.39925,0.3559689952783316,0.3559689952783316,3,1,12.62954966666667,37.888649,0.2986690347680022,0.2986690347680022,224.80282667925,2023.22544011325,0.1063169162271148,0.1063169162271149),
)

$startRow = 2
$r = $startRow
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
